{"js": "// Replace the cover-page placeholder \"NOMES DOS ALUNOS DO GRUPO\" with the\n// first student's name (\"MARIA FERNANDA GALDINO\") and add a second student\n// name (\"CAIO MARINHO DO REIS\") as a new paragraph right after it, copying\n// the same paragraph/run formatting (centered, bold, Arial east-asia font).\n\nconst body = context.document.body;\n\n// Locate the paragraph that currently holds the placeholder text.\nconst searchResults = body.search(\"NOMES DOS ALUNOS DO GRUPO\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"NOMES DOS ALUNOS DO GRUPO\" in the document body.');\n}\n\nconst targetRange = searchResults.items[0];\nconst paragraphs = targetRange.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst namesParagraph = paragraphs.items[0];\n\n// Insert the second student's name as a new paragraph right after the\n// placeholder paragraph; it inherits the placeholder paragraph's formatting\n// (centered, bold, Arial east-asia font, szCs 24).\nconst secondStudentParagraph = namesParagraph.insertParagraph(\"CAIO MARINHO DO REIS\", Word.InsertLocation.after);\nsecondStudentParagraph.alignment = Word.Alignment.centered;\nsecondStudentParagraph.font.bold = true;\n\n// Replace the placeholder text with the first student's name, keeping the\n// run/paragraph formatting intact.\ntargetRange.insertText(\"MARIA FERNANDA GALDINO\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Replace the cover-page placeholder \"NOMES DOS ALUNOS DO GRUPO\" with the\n# first student's name (\"MARIA FERNANDA GALDINO\") and add a second student\n# name (\"CAIO MARINHO DO REIS\") as a new paragraph right after it, copying\n# the same paragraph/run formatting (centered, bold, Arial east-asia font).\n\n$d = $word.ActiveDocument\n\n# Locate the placeholder text and replace it with the first student's name.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"NOMES DOS ALUNOS DO GRUPO\")\nif (-not $found) {\n    throw 'Could not find \"NOMES DOS ALUNOS DO GRUPO\" in the document.'\n}\n$findRange.Text = \"MARIA FERNANDA GALDINO\"\n\n# Re-locate that paragraph by index (COM Paragraphs collection is 1-based)\n# now that the text has changed.\n$targetIndex = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*MARIA FERNANDA GALDINO*\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq 0) {\n    throw \"Could not locate the updated paragraph.\"\n}\n\n# Insert a new paragraph right after it, then fill it with the second\n# student's name; the new paragraph inherits the source paragraph/run\n# formatting (centered, bold, Arial east-asia font, szCs 24).\n$namesParagraph = $d.Paragraphs.Item($targetIndex)\n$insertRange = $namesParagraph.Range\n$insertRange.Collapse(0)\n$insertRange.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Item($targetIndex + 1)\n$newRange = $newParagraph.Range\n$newRange.InsertBefore(\"CAIO MARINHO DO REIS\")\n"}
